# "Fruta / hortaliza, semanal" -- weekly refresh of the price records.
# The underlying data rows (2-18) are re-shuffled: every row's final content
# is an exact copy of some other row's original content (a permutation),
# i.e. this models the weekly roll of the dataset where each record slides
# into a different slot. We snapshot every source row first (so reads never
# see a partially-written destination), then write the snapshots out to
# their destination rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- snapshot every row (columns A:R) before writing anything ---
$row2  = $ws.Range("A2:R2").Value2
$row3  = $ws.Range("A3:R3").Value2
$row4  = $ws.Range("A4:R4").Value2
$row5  = $ws.Range("A5:R5").Value2
$row6  = $ws.Range("A6:R6").Value2
$row7  = $ws.Range("A7:R7").Value2
$row8  = $ws.Range("A8:R8").Value2
$row9  = $ws.Range("A9:R9").Value2
$row10 = $ws.Range("A10:R10").Value2
$row11 = $ws.Range("A11:R11").Value2
$row12 = $ws.Range("A12:R12").Value2
$row13 = $ws.Range("A13:R13").Value2
$row14 = $ws.Range("A14:R14").Value2
$row15 = $ws.Range("A15:R15").Value2
$row16 = $ws.Range("A16:R16").Value2
$row17 = $ws.Range("A17:R17").Value2
$row18 = $ws.Range("A18:R18").Value2

# --- write each snapshot to its new (destination) row ---
$ws.Range("A2:R2").Value2   = $row10
$ws.Range("A3:R3").Value2   = $row9
$ws.Range("A4:R4").Value2   = $row6
$ws.Range("A5:R5").Value2   = $row4
$ws.Range("A6:R6").Value2   = $row17
$ws.Range("A7:R7").Value2   = $row8
$ws.Range("A8:R8").Value2   = $row7
$ws.Range("A9:R9").Value2   = $row15
$ws.Range("A10:R10").Value2 = $row13
$ws.Range("A11:R11").Value2 = $row14
$ws.Range("A12:R12").Value2 = $row12
$ws.Range("A13:R13").Value2 = $row11
$ws.Range("A14:R14").Value2 = $row18
$ws.Range("A15:R15").Value2 = $row2
$ws.Range("A16:R16").Value2 = $row3
$ws.Range("A17:R17").Value2 = $row16
$ws.Range("A18:R18").Value2 = $row5
